$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 23 (shifts old rows 23-29 down to 24-30),
# inheriting formatting from the surrounding cells.
$ws.Rows(23).Insert()

# Row 2: station changed to USGS gage 073802332 (offset stays 0)
$ws.Range("A2").Value = "USGS"
$ws.Range("B2").Value = "073802332"

# Row 7: station id changed to 82740 (now stored as text) with a new offset
$ws.Range("B7").Value = "82740"
$ws.Range("C7").Value = -0.16

# Row 11: offset updated
$ws.Range("C11").Value = -4.5

# Newly inserted row 23: USACE station 76305 with its offset
$ws.Range("A23").Value = "USACE"
$ws.Range("B23").Value = "76305"
$ws.Range("C23").Value = -0.08

# Update the active selection to match the latest edit location
[void]$ws.Range("D23").Select()
